# edit.ps1 - applies the diff:
#  1) Paragraph 1: split the opening run so that "Voici la synthèse de la
#     Discussion précédente : creat-organisation artefact - 3" becomes
#     "Voici la synthèse de la Discussion précédente :" + " .." (the latter
#     carrying an explicit fr-FR language mark), keeping the trailing
#     ". Continuons. Prends en compte que : " run untouched.
#  2) The "Les Cahiers de charge ..." paragraph is rewritten and merged
#     with the following "Considère tous les détails." paragraph.

$d = $word.ActiveDocument

# --- Change 1: first paragraph -------------------------------------------------
$p1 = $d.Paragraphs(1)
$rng1 = $d.Range($p1.Range.Start, $p1.Range.End)

$xml1 = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p w:rsidR="00757A0E" w:rsidRDefault="00757A0E" w:rsidP="00757A0E"><w:pPr><w:pStyle w:val="NormalWeb"/></w:pPr><w:r><w:t>Voici la synthèse de la Discussion précédente :</w:t></w:r><w:r><w:rPr><w:lang w:val="fr-FR"/></w:rPr><w:t xml:space="preserve"> ..</w:t></w:r><w:r><w:rPr><w:rStyle w:val="lev"/></w:rPr><w:t xml:space="preserve">. Continuons. Prends en compte que : </w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

[void]$rng1.InsertXML($xml1)

# --- Change 2: "Les Cahiers de charge..." paragraph merged with "Considère" ---
# Locate the paragraphs by their text so indices stay correct even though the
# first change may shift paragraph numbering internally.
$targetStart = $null
$targetEnd = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $txt = $d.Paragraphs($i).Range.Text
    if ($txt.StartsWith("Les Cahiers de charge complets")) {
        $targetStart = $d.Paragraphs($i).Range.Start
        $targetEnd = $d.Paragraphs($i + 1).Range.End
    }
}

$rng2 = $d.Range($targetStart, $targetEnd)

$xml2 = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p w:rsidR="00757A0E" w:rsidRDefault="00757A0E" w:rsidP="00757A0E"><w:pPr><w:pStyle w:val="NormalWeb"/></w:pPr><w:r><w:t xml:space="preserve">Les Cahiers de charge complets, la structuration du projet, et certaines synthèses des discussions précédentes </w:t></w:r><w:r><w:rPr><w:lang w:val="fr-FR"/></w:rPr><w:t xml:space="preserve">ont été </w:t></w:r><w:r><w:t>partagés</w:t></w:r><w:r><w:rPr><w:lang w:val="fr-FR"/></w:rPr><w:t xml:space="preserve">. </w:t></w:r><w:r><w:t>Considère tous les détails.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

[void]$rng2.InsertXML($xml2)
